$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(4)
$tbl = $sh.Table
$cell = $tbl.Cell(2, 2)
$tr = $cell.Shape.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf(",Morocco(ANRT),Philippines(NTC)")
$start = $idx + 1
$len = ",Morocco(ANRT),Philippines(NTC)".Length
$sub = $tr.Characters($start, $len)
try {
    $sub.Font.Strikethrough = 1
    Write-Host "set ok: $($sub.Font.Strikethrough)"
} catch { Write-Host "ERR set Strikethrough: $_" }
